$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows describing World / WorldStage related packets.
# Order of cell writes matters for shared-string table ordering:
# WorldList, WorldStageList, LIST:WorldStagePacket, LIST:WorldPacket
$ws.Range("A31").Value = "WorldList"
$ws.Range("A32").Value = "WorldStageList"
$ws.Range("B32").Value = "LIST:WorldStagePacket"
$ws.Range("B31").Value = "LIST:WorldPacket"
$ws.Range("C31").Value = "new()"
$ws.Range("C32").Value = "new()"
$ws.Range("E31").Value = "Packet"
$ws.Range("E32").Value = "Packet"

$ws.Range("B29").Select()
